$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...heading as well as the stated missing parts." -> "...requested
#    missing parts."
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "stated missing parts", $true, $false, $false, $false, $false,
    $true, 1, $false, "requested missing parts", 2)

# ---------------------------------------------------------------------
# 2) Quoted paragraph: "...provide the users with traditional 3D view..."
#    -> "...provide the users only with traditional 3D view..."
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "provide the users with traditional 3D view",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "provide the users only with traditional 3D view", 2)

# ---------------------------------------------------------------------
# 3) Move the hidden "_GoBack" bookmark from the "...three most typical."
#    sentence to the empty paragraph right after "...best plus a simple
#    tutorial on its usage." Bookmark names are unique, so adding it in
#    the new spot removes it from the old spot automatically.
# ---------------------------------------------------------------------
$goBackTarget = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -eq "`r") {
        $prev = $para.Previous()
        if ($prev -ne $null -and $prev.Range.Text -like "*best plus a simple tutorial on its usage.*") {
            $goBackTarget = $para
            break
        }
    }
}
if ($goBackTarget -eq $null) {
    # Fallback: known paragraph index in the original document.
    $goBackTarget = $d.Paragraphs(60)
}
$d.Bookmarks.Add("_GoBack", $goBackTarget.Range)

# ---------------------------------------------------------------------
# 4) "...examples of three most typical ." -> "...three most typical
#    protein-protein interaction modes." Rebuilding the whole paragraph
#    also clears the (now orphaned) gramStart/gramEnd proofing markers
#    that used to wrap "typical" / trail the sentence.
# ---------------------------------------------------------------------
$typicalPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*three most typical*") {
        $typicalPara = $para
        break
    }
}
$pStart = $typicalPara.Range.Start
$pEnd = $typicalPara.Range.End
$fullRange = $d.Range($pStart, $pEnd)
$fullRange.Delete()
$insertionPoint = $d.Range($pStart, $pStart)
$insertionPoint.InsertAfter(
    "We were trying to demonstrate the usage of the tool in the Results " + `
    "and Discussion section where we showcase examples of three most " + `
    "typical protein-protein interaction modes.`r")

# ---------------------------------------------------------------------
# 5) Insert " like Residue Matrix and Contact Zone Graph" right after
#    "...additional specialized views" (before the following period).
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "additional specialized views. Our tool is based",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "additional specialized views like Residue Matrix and Contact Zone Graph. Our tool is based",
    2)

# ---------------------------------------------------------------------
# 6) Merge the split "K" / "ozlíková" runs of the reference author's
#    name back into a single "Kozlíková" run (no text change).
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Kozlíková", $true, $false, $false, $false, $false,
    $true, 1, $false, "Kozlíková", 2)
